$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.431.89"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").Value = "3.577.14"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("E7").Value = "  +11.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.410"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("E9").Value = "  +6.83%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "3.575.50"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").Value = "4.241.71"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").Value = "96.271.04"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000260"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "

$ws.Range("D18").Value = "3.561.08"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.518"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "502.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.54%  "

$ws.Range("D29").Value = "3.768.89"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.154"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("E34").Value = "  +3.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "626.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.566"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.26%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.905"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("E44").Value = "  +5.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.60%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.04%  "
